# Update seed-data passwords in the "Senders" test data sheet.
# The workbook's sole worksheet holds ID / Password / Username columns
# (A/B/C) across rows 2-11. We rewrite the Password column (B) with
# strengthened values, in the same order the original author produced
# them in (B4, B3, B2, then B5..B11 top-to-bottom) so the shared-string
# table layout matches the saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "Blynas1!!"
$ws.Range("B3").Value = "Kiaušinis2!?"
$ws.Range("B2").Value = "Viščiukas.3"
$ws.Range("B5").Value = "Pasas13:"
$ws.Range("B6").Value = "Pasuotikamuolį??4"
$ws.Range("B7").Value = "Krepšininkas17ĄĄ?"
$ws.Range("B8").Value = "SkėtyjeSkylė44?"
$ws.Range("B9").Value = "LietusSuSaule!!3"
$ws.Range("B10").Value = "SaulėSuLietumi?OI4"
$ws.Range("B11").Value = "Lietuviųkalba55?"

# Move the active selection to match the state the file was saved in.
$ws.Range("C10").Select() | Out-Null
